$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handback" — refresh the localization-status report:
#   * Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   * The Latest Handback DateTime is refreshed to the new handback timestamp
#   * The stale-handback-version error notes are cleared (handback is now current)
#   * The Status / Error Detail columns are resized to fit the new text
# ---------------------------------------------------------------------------

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns ---------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-13 02:43:14"
$wsZhCn.Range("K3").Value = "2016-08-13 02:43:14"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-13 02:43:24"
$wsDeDe.Range("K3").Value = "2016-08-13 02:43:24"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("P3").Value = ""

# --- Column widths: Status columns grow, Error Detail columns shrink --------
$statusColumnWidth = 29.166666666666668
$errorColumnWidth  = 12.833333333333334

$wsOverview.Columns.Item(5).ColumnWidth = $statusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColumnWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $statusColumnWidth
$wsZhCn.Columns.Item(16).ColumnWidth = $errorColumnWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $statusColumnWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $errorColumnWidth
